# Weekly update: insert a new data row for the week of 2023-03-31
# (Excel serial date 45016) ahead of the existing series, shifting the
# remaining rows (old 209..308) down to (210..309).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 209 - this pushes the old row 209 (and
# everything below it) down by one row, preserving all of their values.
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new week's record.
$ws.Cells.Item(209, 1).Value = 8
$ws.Cells.Item(209, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(209, 3).Value = "Coquimbo"
$ws.Cells.Item(209, 4).Value = 45016
$ws.Cells.Item(209, 5).Value = 4
$ws.Cells.Item(209, 6).Value = 100112037
$ws.Cells.Item(209, 7).Value = "Cebollín"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 1300
$ws.Cells.Item(209, 11).Value = 1000
$ws.Cells.Item(209, 12).Value = 1200
$ws.Cells.Item(209, 13).Value = 1100
$ws.Cells.Item(209, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(209, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(209, 16).Value = 183
$ws.Cells.Item(209, 17).Value = 6
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Make sure the new date cell carries the same date style as the rest of
# column D (row insert should already have copied formatting, but be
# explicit/defensive about it).
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
